$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# 1. Update activation date
Replace-Text "Ativação: 01/01/2022" "Ativação: 01/01/2025"

# 2. Update English objectives paragraph
Replace-Text "To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering." "Introduce students to Biochemical Engineering, the characteristics of the profession, and guide them regarding the responsibilities and areas of practice of a Biochemical Engineer. Additionally, develop in students a macro view of the types and stages of an industrial bioprocess, and finally, guide them on the role of the Biochemical Engineer in industry, research and teaching, as well as entrepreneurship and innovation in engineering."

# 3. Fix spacing in Portuguese summary ("Programa resumido")
Replace-Text "1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia4. Áreas de atuação do Engenheiro Bioquímico5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos)8. Visita supervisionada." "1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia 4. Áreas de atuação do Engenheiro Bioquímico 5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos) 8. Visita supervisionada."

# 4. Rewrite English summary ("Programa resumido")
Replace-Text "1. History of Biochemical Engineering 2. Biochemical Engineering: Definitions and Concepts 3. Engineering Job Market 4. Areas of practice of the Biochemical Engineer 5. The Bioprocess Industry 6. Production Scales 7. Case studies (biotechnological processes) 8. Supervised visit." "1. History of Biochemical Engineering; 2. Biochemical Engineering: Definitions and concepts; 3. Job market for Engineering; 4. Areas of practice for the Biochemical Engineer; 5. The Bioprocess Industry; 6. Production scalesCase studies (biotechnological processes); 7. Supervised visit."

# 5. Fix spacing & duplicate the Portuguese full program ("Programa")
Replace-Text "1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso." "1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil 3. Atribuições e áreas de atuação do Engenheiro Bioquímico 4. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos). 8. Empreendedorismo e Inovação em Engenharia. 9. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso."

# 6. Rewrite English full program ("Programa")
Replace-Text "1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts – enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry – types of industries, equipment, installations, main unit operations6.Production scales – laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation – visits to laboratories and bioprocess industry" "1.History of Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market for Engineering in Brazil3.Responsibilities and areas of practice for the Biochemical Engineer4.Definitions and concepts – enzymatic process, generic fermentation process, transforming agents, bioreactor, raw materials, types of substrates, substrate-to-product conversion, types of biotechnological products, product recovery, among others.5.The Bioprocess Industry – types of industries, equipment, facilities, main unit operations.6.Production scales – laboratory, pilot, industrial.7.Case studies (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visits – visits to laboratories and bioprocess industry."
